# The edit permutes the content of the data rows (rows 3-15) on the
# "Artfynd" sheet: every data row's identity (species record) moves to a
# different row position. Most columns (C, I, P, S, T, U, V, W, Y, Z, AA,
# AB, AD, AE, AG, AT, AW, AX, AY) are identical across all of rows 3-15,
# so only the columns that actually carry distinguishing data need to be
# moved: A, B, D, E, F, G, H, Q, R (the record's fields) and AC (a free
# text comment that is only present on one row at a time).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow = 15

# Columns (by number) whose values differ row to row and therefore need
# to be carried along with the permutation.
$movingCols = @(1, 2, 4, 5, 6, 7, 8, 17, 18)   # A,B,D,E,F,G,H,Q,R
$commentCol = 29                                # AC

# Mapping of source row -> destination row for the permutation (derived
# from matching each record's unique Id, column A, between the before
# and after states).
$rowMap = @{
    3  = 12
    4  = 7
    5  = 11
    6  = 13
    7  = 8
    8  = 4
    9  = 3
    10 = 9
    11 = 5
    12 = 10
    13 = 14
    14 = 15
    15 = 6
}

# 1) Snapshot the current ("before") values of the moving columns, plus
#    the comment column, for every affected row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $movingCols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot["$r-$commentCol"] = $ws.Cells.Item($r, $commentCol).Value()
}

# 2) Clear out the comment column across the affected rows; it will be
#    re-populated only where the snapshot says it should land.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $commentCol).ClearContents()
}

# 3) Write the snapshotted values to their destination rows.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $destRow = $rowMap[$r]

    foreach ($c in $movingCols) {
        $ws.Cells.Item($destRow, $c).Value = $snapshot["$r-$c"]
    }

    $commentVal = $snapshot["$r-$commentCol"]
    if ($commentVal -ne $null -and $commentVal -ne "") {
        $ws.Cells.Item($destRow, $commentCol).Value = $commentVal
    }
}
